$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "Play Always Hot Deluxe Slot Game for Free - Review" "Play Always Hot Deluxe Free: Classic Slot Machine Review"
Replace-Text "Simple and straightforward gameplay" "Classic slot machine with simple and straightforward gameplay"
Replace-Text "Non-progressive jackpot of up to 60,000 coins" "Offers a non-progressive jackpot of up to 60,000 coins"
Replace-Text "Gambling feature for added excitement" "Gamble feature adds an element of risk and excitement"
Replace-Text "Suitable for players who enjoy low or high stakes" "Suitable for players who enjoy low or high stakes betting"
Replace-Text "Limited features compared to other slot games" "Limited features compared to other slot machine games"
Replace-Text "No free spins or bonus rounds" "May not appeal to players looking for more complex gameplay"
Replace-Text "Read our unbiased review of Always Hot Deluxe, a classic slot machine with a non-progressive jackpot of up to 60,000 coins. Play for free and learn more!" "Play Always Hot Deluxe for free and enjoy classic slot machine gameplay with a chance to win big."
